# Fruta / hortaliza, semanal
# Insert a new weekly record at row 266 for "Granada" (Vega Modelo de Temuco),
# shifting existing rows 266-320 down to 267-321.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 266, pushing the rest of the table down.
$ws.Range("A266:T266").EntireRow.Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A266").Value = 10
$ws.Range("B266").Value = "Vega Modelo de Temuco"
$ws.Range("C266").Value = "La Araucanía"
$ws.Range("D266").Value = 45258
$ws.Range("D266").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E266").Value = 9
$ws.Range("F266").Value = "Fruta"
$ws.Range("G266").Value = 100104
$ws.Range("H266").Value = "Frutos de pepita"
$ws.Range("I266").Value = 100104001
$ws.Range("J266").Value = "Granada"
$ws.Range("K266").Value = "Wonderfull"
$ws.Range("L266").Value = "Primera"
$ws.Range("M266").Value = 60
$ws.Range("N266").Value = 15000
$ws.Range("O266").Value = 15000
$ws.Range("P266").Value = 15000
$ws.Range("Q266").Value = "`$/bandeja 10 kilos granel"
$ws.Range("R266").Value = "Provincia de Limarí"
$ws.Range("S266").Value = 1500
$ws.Range("T266").Value = 10
